$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix C27: was stored as a text value, should be a real number.
$ws.Range("C27").Value = 58475425000137

# Append the new log row (row 28).
$ws.Range("A28").Value = "03/07/2025 17:08:13"
$ws.Range("B28").Value = "Ima Industria"

# CNPJ has a leading zero, so it must stay text - prefix with an apostrophe
# (classic Excel "force text" entry) otherwise Excel would coerce the
# numeric-looking string to a number and the leading zero would be lost.
# Reset the style afterwards so the quote-prefix formatting Excel applies
# doesn't leave a stray style reference on the cell.
$ws.Range("C28").Value = "'04252502000160"
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").Value = "denissonfhsilva@gmail.com"
$ws.Range("E28").Value = "893-ExtratoMensal-052025.pdf"
